# Update "想去人数" (F) counts (and, for two rows, the "最低票价" (G)
# column which flipped from a numeric price to "不可售") across the
# three data sheets that carry this event list: "展览", "演出" and the
# combined "全部类型" view. "本地生活" has no data rows and is untouched.

$wb = $excel.ActiveWorkbook

# ---- 展览 ("Exhibitions") ----
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value  = 3373
$ws.Cells.Item(4, 6).Value  = 2440
$ws.Cells.Item(6, 6).Value  = 340
$ws.Cells.Item(7, 6).Value  = 1387
$ws.Cells.Item(8, 6).Value  = 1095
$ws.Cells.Item(9, 6).Value  = 302
$ws.Cells.Item(10, 6).Value = 515
$ws.Cells.Item(13, 6).Value = 100
$ws.Cells.Item(15, 6).Value = 8582
$ws.Cells.Item(16, 6).Value = 373
$ws.Cells.Item(17, 6).Value = 2483
$ws.Cells.Item(22, 6).Value = 585
$ws.Cells.Item(26, 6).Value = 2016
$ws.Cells.Item(26, 7).Value = "不可售"
$ws.Cells.Item(27, 6).Value = 2055
$ws.Cells.Item(29, 6).Value = 1746
$ws.Cells.Item(33, 6).Value = 28
$ws.Cells.Item(34, 6).Value = 43
$ws.Cells.Item(35, 6).Value = 87
$ws.Cells.Item(36, 6).Value = 182
$ws.Cells.Item(38, 6).Value = 301
$ws.Cells.Item(40, 6).Value = 237
$ws.Cells.Item(41, 6).Value = 410
$ws.Cells.Item(42, 6).Value = 151
$ws.Cells.Item(44, 6).Value = 259

# ---- 演出 ("Performances") ----
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 16

# ---- 全部类型 ("All types" - combined view) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value  = 3373
$ws.Cells.Item(4, 6).Value  = 2440
$ws.Cells.Item(6, 6).Value  = 340
$ws.Cells.Item(7, 6).Value  = 1387
$ws.Cells.Item(9, 6).Value  = 1095
$ws.Cells.Item(10, 6).Value = 302
$ws.Cells.Item(11, 6).Value = 515
$ws.Cells.Item(13, 6).Value = 100
$ws.Cells.Item(15, 6).Value = 8582
$ws.Cells.Item(16, 6).Value = 373
$ws.Cells.Item(17, 6).Value = 2483
$ws.Cells.Item(18, 6).Value = 16
$ws.Cells.Item(23, 6).Value = 585
$ws.Cells.Item(27, 6).Value = 2016
$ws.Cells.Item(27, 7).Value = "不可售"
$ws.Cells.Item(28, 6).Value = 2055
$ws.Cells.Item(29, 6).Value = 1746
$ws.Cells.Item(33, 6).Value = 28
$ws.Cells.Item(34, 6).Value = 43
$ws.Cells.Item(35, 6).Value = 87
$ws.Cells.Item(36, 6).Value = 182
$ws.Cells.Item(38, 6).Value = 301
$ws.Cells.Item(40, 6).Value = 237
$ws.Cells.Item(41, 6).Value = 410
$ws.Cells.Item(46, 6).Value = 151
$ws.Cells.Item(49, 6).Value = 259
